$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "kNN"
$ws.Range("B2").Value = 0.84
$ws.Range("C2").Value = 58.08
$ws.Range("D2").Value = 6760.02
$ws.Range("E2").Value = 82.21

$ws.Range("D3").Select()
